# Update the "想去人数" (interest/attendance count) figures in column F
# for the affected rows on sheets "展览" (sheet1), "演出" (sheet2) and
# "全部类型" (sheet4), matching the regenerated data from the upstream
# bilibili feed (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7754
$ws1.Range("F6").Value = 134
$ws1.Range("F23").Value = 124
$ws1.Range("F27").Value = 2578
$ws1.Range("F38").Value = 4780
$ws1.Range("F43").Value = 933
$ws1.Range("F44").Value = 327

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 7
$ws2.Range("F16").Value = 114

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 7754
$ws4.Range("F7").Value = 134
$ws4.Range("F9").Value = 7
$ws4.Range("F27").Value = 124
$ws4.Range("F29").Value = 2578
$ws4.Range("F38").Value = 4780
$ws4.Range("F45").Value = 933
$ws4.Range("F46").Value = 327
